# The document's "References" section ends with a long auto-generated
# BIBLIOGRAPHY field/content-control, immediately followed by a run of
# individual CITATION field content controls - one per paragraph, each
# paragraph's text being just "(Author, Year)". A "_GoBack" bookmark
# pair used to sit in the middle of that run of citation paragraphs.
#
# This commit deletes every one of those trailing citation paragraphs
# (the whole run at the end of the document), which leaves the
# "_GoBack" bookmark pair as the only thing still sitting between the
# end of the bibliography content control and the section break.

$d = $word.ActiveDocument

# Walk backwards from the last paragraph while paragraphs look like a
# bare "(...)" citation render (no other visible text), to robustly
# find where that trailing run of citation paragraphs begins -
# regardless of exactly how many of them there are.
$count = $d.Paragraphs.Count
$firstCitationIndex = $count + 1
for ($i = $count; $i -ge 1; $i--) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text.Trim()
    $looksLikeCitation = $text.StartsWith("(") -and $text.EndsWith(")")
    if (-not $looksLikeCitation) {
        break
    }
    $firstCitationIndex = $i
}

if ($firstCitationIndex -le $count) {
    $startOfDeletion = $d.Paragraphs.Item($firstCitationIndex).Range.Start
    $endOfDeletion = $d.Content.End
    $deleteRange = $d.Range($startOfDeletion, $endOfDeletion)
    $deleteRange.Delete()
}
